# Weekly update: insert a new data row (row 3) for Damasco / Vega Modelo de
# Temuco, shifting the existing rows 3-47 down to 4-48.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 3 -- this pushes rows 3..47 down
# to 4..48 and extends the used range to A1:T48.
$ws.Rows.Item(3).Insert()

# Populate the newly-inserted row 3 with this week's entry.
$ws.Range("A3").Value = 10
$ws.Range("B3").Value = "Vega Modelo de Temuco"
$ws.Range("C3").Value = "La Araucanía"
$ws.Range("D3").Value = 44552
$ws.Range("E3").Value = 9
$ws.Range("F3").Value = "Fruta"
$ws.Range("G3").Value = 100103
$ws.Range("H3").Value = "Frutos de hueso (carozo)"
$ws.Range("I3").Value = 100103003
$ws.Range("J3").Value = "Damasco"
$ws.Range("K3").Value = "Dina"
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 55
$ws.Range("N3").Value = 20000
$ws.Range("O3").Value = 22000
$ws.Range("P3").Value = 21091
$ws.Range("Q3").Value = "$/caja 18 kilos"
$ws.Range("R3").Value = "Provincia de Quillota"
$ws.Range("S3").Value = 1172
$ws.Range("T3").Value = 18
